$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.754.34'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '2.541.46'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.14'
$ws.Range('E5').Value = '  +4.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.35'
$ws.Range('E6').Value = '  -3.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.578'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.42'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.71'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.113'
$ws.Range('E13').Value = '  -0.84%  '
$ws.Range('D14').Value = '2.932.73'
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.87'
$ws.Range('E15').Value = '  +4.61%  '
$ws.Range('D16').Value = '2.543.72'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.868'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').Value = '42.824.47'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.15'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('D21').Value = '0.0₃0969'
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.10'
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '252.73'
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.02'
$ws.Range('E25').Value = '  -2.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.12'
$ws.Range('E26').Value = '  -2.27%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  +3.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.03'
$ws.Range('E29').Value = '  +3.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.26'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '155.74'
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('E33').Value = '  +0.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.33'
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.08'
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0791'
$ws.Range('E36').Value = '  -1.04%  '
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('E38').Value = '  -2.94%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.119'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.62'
$ws.Range('E40').Value = '  -4.93%  '
$ws.Range('B41').Value = 'ApeXProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.36'
$ws.Range('E41').Value = '  +11.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.86'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.36'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0302'
$ws.Range('E45').Value = '  -1.05%  '
$ws.Range('D46').Value = '2.025.61'
$ws.Range('E46').Value = '  -2.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '84.73'
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.93'
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('D49').Value = '2.788.83'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.33'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('E51').Value = '  -0.74%  '
